# dlgTruckInsurance_pagProductData.xlsx
# "msz - repplacing AppiumLibrary by native python-appium-client (part 1)"
#
# The hard-coded start-date test value ("05/01/2025") is replaced by a
# dynamic placeholder token that the new python-appium-client based test
# framework resolves at run time (32 days from today).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Replace the literal start date in B5 with the dynamic <TODAY ...> token.
$ws.Range("B5").Value = "<TODAY +32,+0,+0,'%m/%d/%Y'>"

# 2) Column B (best-fit) widens to accommodate the longer placeholder text.
$ws.Columns.Item(2).ColumnWidth = 27.33

# 3) The second screenshot picture shrinks a little so the page keeps fitting
#    next to the now-wider column B.
$sh = $ws.Shapes.Item(2)
$sh.Width = 729.6800639763779

# 4) Leave the selection where the editor's cursor ended up.
$ws.Range("E8").Select() | Out-Null
